$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Invalid_SignUp _Data"

# Row 2: fill in numeric values (previously empty but styled cells)
$ws.Cells.Item(2, 1).Value = 123
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 123

# Row 3: new negative-scenario data row with a hyperlink in column A
$ws.Cells.Item(3, 2).Value = "Jis"
$ws.Cells.Item(3, 3).Value = "Jay"
$ws.Cells.Item(3, 4).Value = "abcd"

$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:testnum123@yopmail.com", "", "", "testnum123@yopmail.com")

# Update selection to match final state
$ws.Range("C15").Select()
